# edit.ps1 - applies the resume edit described by the diff:
#  1. Removes one of the two blank paragraphs before "EDUCATION"
#  2. Removes the blank (ind=720) paragraph before "TECHNICAL SKILLS"
#  3. Reshuffles the Technical Skills bullet list:
#       Bootstrap is dropped, JavaScript/jQuery.../Node.js each shift up one
#       slot, and a new "Git" bullet is appended; the _GoBack bookmark moves
#       to sit right after "CSS"
#  4. Replaces the "Main page for Maxx Potential company:" freelance entry
#     with a new "Word typing game:" entry (new netlify link + a new
#     "source code:" link to the GitHub repo)

$d = $word.ActiveDocument

function Get-ParaBounds($para) {
    return ,@($para.Range.Start, $para.Range.End)
}

function Set-ParaText($para, [string]$text) {
    $s = $para.Range.Start
    $e = $para.Range.End
    $r = $d.Range($s, $e - 1)
    $r.Text = $text
    return $r
}

# Insert a bookmark immediately at $pos (a position that sits on a paragraph
# boundary) without it snapping to the nearest paragraph edge: type a
# placeholder character right after $pos, drop the (now mid-paragraph,
# non-boundary) bookmark in the gap, then remove the placeholder again.
function Add-BookmarkAt($pos, [string]$name) {
    $tmp = $d.Range($pos, $pos)
    $tmp.InsertAfter("X")
    $bm = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $bm) | Out-Null
    $del = $d.Range($pos, $pos + 1)
    $del.Delete()
}

# Insert a hyperlink immediately at $pos the same boundary-safe way: type a
# placeholder character, wrap the hyperlink around that one character, then
# the placeholder text becomes the (only) display text so nothing needs
# deleting afterwards.
function Add-HyperlinkAt($pos, [string]$address, [string]$displayText) {
    $tmp = $d.Range($pos, $pos)
    $tmp.InsertAfter("X")
    $linkRng = $d.Range($pos, $pos + 1)
    $h = $d.Hyperlinks.Add($linkRng, $address, $null, $null, $displayText)
    return $h
}

# ---------------------------------------------------------------------
# 1) Remove one of the two blank paragraphs right before "EDUCATION"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("EDUCATION", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$eduPara = $rng.Paragraphs(1)
$blankPara = $eduPara.Previous()
$blankPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) Remove the blank (ind=720) paragraph right before "TECHNICAL SKILLS"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("TECHNICAL SKILLS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tsPara = $rng.Paragraphs(1)
$blankPara2 = $tsPara.Previous()
$blankPara2.Range.Delete()

# ---------------------------------------------------------------------
# 3) Technical skills bullet-list reshuffle
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("TECHNICAL SKILLS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$heading = $rng.Paragraphs(1)

$pCpp = $heading.Next()
$pHtml = $pCpp.Next()
$pCss = $pHtml.Next()
$pBootstrap = $pCss.Next()
$pJs = $pBootstrap.Next()
$pJquery = $pJs.Next()
$pNode = $pJquery.Next()

# Shift the text of each bullet up one slot (Bootstrap disappears, Git appended)
Set-ParaText $pNode "Git" | Out-Null
Set-ParaText $pJquery "Node.js" | Out-Null
Set-ParaText $pJs "jQuery, GSAP, React.js" | Out-Null
Set-ParaText $pBootstrap "JavaScript " | Out-Null

# Move the _GoBack bookmark to sit right after the CSS text
$cssEnd = $pCss.Range.End
Add-BookmarkAt ($cssEnd - 1) "_GoBack"

# ---------------------------------------------------------------------
# 4) Freelance projects: "Main page for Maxx..." -> "Word typing game:"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Main page for Maxx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$maxxPara = $rng.Paragraphs(1)
Set-ParaText $maxxPara "Word typing game:" | Out-Null

$linkPara = $maxxPara.Next()
$linkParaStart = $linkPara.Range.Start
$linkParaEnd = $linkPara.Range.End

$hIndex = -1
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks($i)
    if ($h.Range.Start -ge $linkParaStart -and $h.Range.End -le $linkParaEnd) {
        $hIndex = $i
    }
}
$codepenLink = $d.Hyperlinks($hIndex)
$hStart = $codepenLink.Range.Start
if ($hStart -gt $linkParaStart) {
    $spacesRng = $d.Range($linkParaStart, $hStart)
    $spacesRng.Delete()
}
$codepenLink.Address = "https://artwordtyping.netlify.com/"
$codepenLink.TextToDisplay = "https://artwordtyping.netlify.com/"
$linkPara.LeftIndent = 36

# Insert a new "source code: <link>" paragraph right after it
$linkPara.Range.InsertParagraphAfter()
$srcPara = $linkPara.Next()
$srcPara.Range.Text = "source code: X"
$srcPara.LeftIndent = 36
$srcLinkPos = $srcPara.Range.End - 2
Add-HyperlinkAt $srcLinkPos "https://github.com/earthddx/wordTypingGame" "https://github.com/earthddx/wordTypingGame" | Out-Null

Write-Output "done"
